$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Order matters: a couple of the new values coincide with old values that
# are themselves replaced later in the document, so those "source" cells
# must be handled first to avoid the later rule re-matching freshly
# inserted text (e.g. 63÷2= is both an original value and the result of
# 87÷6=, so the original 63÷2= must become 38÷2= before 87÷6= becomes 63÷2=).
Replace-Exact "96÷2=" "95÷5="
Replace-Exact "14÷2=" "64÷6="
Replace-Exact "36÷3=" "65÷4="
Replace-Exact "84÷4=" "74÷3="
Replace-Exact "31÷2=" "56÷2="
Replace-Exact "63÷2=" "38÷2="
Replace-Exact "86÷5=" "15÷8="
Replace-Exact "79÷6=" "75÷9="
Replace-Exact "31÷7=" "49÷7="
Replace-Exact "19÷7=" "19÷6="
Replace-Exact "39÷2=" "16÷6="
Replace-Exact "67÷3=" "64÷2="
Replace-Exact "40÷8=" "91÷4="
Replace-Exact "88÷6=" "99÷5="
Replace-Exact "48÷4=" "89÷5="
Replace-Exact "96÷3=" "74÷7="
Replace-Exact "73÷9=" "34÷6="
Replace-Exact "87÷6=" "63÷2="
Replace-Exact "37÷7=" "82÷2="
Replace-Exact "30÷6=" "62÷9="
Replace-Exact "59÷3=" "35÷8="
Replace-Exact "11÷5=" "88÷7="
Replace-Exact "81÷3=" "54÷3="
Replace-Exact "76÷4=" "55÷9="
Replace-Exact "34÷2=" "67÷3="
